# Add a new worksheet "2024-11-27" with a small stock review table,
# mirroring the structure of the existing "2024-11-19"/"2024-11-26" sheets.

$wb = $excel.ActiveWorkbook

# Add the new sheet; Excel inserts it after the active sheet by default,
# but we want it appended at the end (after "2024-11-26").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2024-11-27"

# Header row
$ws.Range("A1").Value = "stock"
$ws.Range("B1").Value = "StartPrice"

# Data rows -- stock codes are kept as text (as in the other sheets),
# prices as numbers. Force text formatting so numeric-looking codes like
# "1231" are not auto-converted to numbers.
$ws.Range("A2:A3").NumberFormat = "@"

$ws.Range("A2").Value = "1231"
$ws.Range("B2").Value = 120

$ws.Range("A3").Value = "1319"
$ws.Range("B3").Value = 114.5

$ws.Range("A1").Select()
